$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Cases by Age Group" - updated case counts by age group
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Cases by Age Group")
$ws1.Cells.Item(2,2).Value  = 267
$ws1.Cells.Item(3,2).Value  = 1269
$ws1.Cells.Item(4,2).Value  = 3446
$ws1.Cells.Item(5,2).Value  = 15119
$ws1.Cells.Item(6,2).Value  = 16730
$ws1.Cells.Item(7,2).Value  = 14622
$ws1.Cells.Item(8,2).Value  = 12266
$ws1.Cells.Item(9,2).Value  = 4407
$ws1.Cells.Item(10,2).Value = 2963
$ws1.Cells.Item(11,2).Value = 1747
$ws1.Cells.Item(12,2).Value = 1139
$ws1.Cells.Item(13,2).Value = 1762
$ws1.Cells.Item(14,2).Value = 13

# ------------------------------------------------------------------
# Sheet 2: "Cases by Gender" - updated case counts by gender
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Cases by Gender")
$ws2.Cells.Item(2,2).Value = 25349
$ws2.Cells.Item(3,2).Value = 49491
$ws2.Cells.Item(4,2).Value = 910

# ------------------------------------------------------------------
# Sheet 3: "Cases by RaceEthnicity" - updated case counts by race/ethnicity
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws3.Cells.Item(2,2).Value = 933
$ws3.Cells.Item(3,2).Value = 12672
$ws3.Cells.Item(4,2).Value = 27728
$ws3.Cells.Item(5,2).Value = 428
$ws3.Cells.Item(6,2).Value = 25581
$ws3.Cells.Item(7,2).Value = 8408

# ------------------------------------------------------------------
# Sheet 4: "Fatalities by Age Group" - updated fatality counts by age group
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Fatalities by Age Group")
$ws4.Cells.Item(3,2).Value  = 10
$ws4.Cells.Item(4,2).Value  = 31
$ws4.Cells.Item(5,2).Value  = 216
$ws4.Cells.Item(6,2).Value  = 700
$ws4.Cells.Item(7,2).Value  = 2075
$ws4.Cells.Item(8,2).Value  = 4796
$ws4.Cells.Item(9,2).Value  = 4042
$ws4.Cells.Item(10,2).Value = 5238
$ws4.Cells.Item(11,2).Value = 5891
$ws4.Cells.Item(12,2).Value = 5913
$ws4.Cells.Item(13,2).Value = 15216

# ------------------------------------------------------------------
# Sheet 5: "Fatalities by Gender" - updated fatality counts by gender
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Fatalities by Gender")
$ws5.Cells.Item(2,2).Value = 18487
$ws5.Cells.Item(3,2).Value = 25646

# ------------------------------------------------------------------
# Sheet 6: "Fatalities by Race-Ethnicity" - updated fatality counts by race/ethnicity
# ------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws6.Cells.Item(2,2).Value = 898
$ws6.Cells.Item(3,2).Value = 4252
$ws6.Cells.Item(4,2).Value = 20453
$ws6.Cells.Item(5,2).Value = 231
$ws6.Cells.Item(6,2).Value = 18276

# ------------------------------------------------------------------
# Update sheet view selections (and which sheet/tab is active).
# Activate in the order that reproduces the final state, finishing
# on the sheet that should end up with tabSelected="1".
# ------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C19").Select()

$ws3.Activate()
$ws3.Range("C17").Select()

$ws4.Activate()
$ws4.Range("C25").Select()

$ws6.Activate()
$ws6.Range("A28:B34").Select()

$ws1.Activate()
$ws1.Range("B20").Select()
